# "Update countries & provincias Spain"
#
# Refreshes the COVID-19 "Pais" sheet with a newer data pull:
#  - timestamp footer bumped to 02:20
#  - a handful of countries got their daily totals updated (Estados Unidos,
#    Japon, Panama, Uruguay, Vietnam) with no re-sort needed
#  - several small blocks of countries were re-sorted (by case count), which
#    shifts which country name/figures land on a given row; a brand-new
#    row for "Eritrea" also appears, pushing a short run of neighboring
#    countries down by one
#
# Values below are written directly against the fixed row/column grid so
# the on-screen result matches the refreshed source exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 31 de Marzo de 2020 a las 02:20'

$ws.Cells.Item(4, 2).Value = 163479
$ws.Cells.Item(4, 3).Value = 19988
$ws.Cells.Item(4, 4).Value = 5506
$ws.Cells.Item(4, 5).Value = 154825

$ws.Cells.Item(22, 1).Value = 'Australia'
$ws.Cells.Item(22, 2).Value = 4460
$ws.Cells.Item(22, 3).Value = 297
$ws.Cells.Item(22, 4).Value = 244
$ws.Cells.Item(22, 5).Value = 4197
$ws.Cells.Item(22, 6).Value = 28
$ws.Cells.Item(22, 7).Value = 2
$ws.Cells.Item(22, 8).Value = 19

$ws.Cells.Item(23, 1).Value = 'Noruega'
$ws.Cells.Item(23, 2).Value = 4445
$ws.Cells.Item(23, 3).Value = 161
$ws.Cells.Item(23, 4).Value = 12
$ws.Cells.Item(23, 5).Value = 4401
$ws.Cells.Item(23, 6).Value = 97
$ws.Cells.Item(23, 7).Value = 6
$ws.Cells.Item(23, 8).Value = 32

$ws.Cells.Item(34, 2).Value = 1953
$ws.Cells.Item(34, 3).Value = 87
$ws.Cells.Item(34, 5).Value = 1473
$ws.Cells.Item(34, 7).Value = 2
$ws.Cells.Item(34, 8).Value = 56

$ws.Cells.Item(46, 4).Value = 9
$ws.Cells.Item(46, 5).Value = 1039
$ws.Cells.Item(46, 6).Value = 43
$ws.Cells.Item(46, 7).Value = 3
$ws.Cells.Item(46, 8).Value = 27

$ws.Cells.Item(79, 4).Value = 25
$ws.Cells.Item(79, 5).Value = 294

$ws.Cells.Item(92, 2).Value = 204
$ws.Cells.Item(92, 3).Value = 10
$ws.Cells.Item(92, 5).Value = 149

$ws.Cells.Item(154, 1).Value = 'Eritrea'
$ws.Cells.Item(154, 3).Value = 3
$ws.Cells.Item(154, 4).Value = 0
$ws.Cells.Item(154, 5).Value = 15

$ws.Cells.Item(155, 1).Value = 'Haiti'
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 1
$ws.Cells.Item(155, 5).Value = 14
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 0

$ws.Cells.Item(156, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(156, 2).Value = 15
$ws.Cells.Item(156, 4).Value = 2
$ws.Cells.Item(156, 5).Value = 12
$ws.Cells.Item(156, 7).Value = 1
$ws.Cells.Item(156, 8).Value = 1

$ws.Cells.Item(157, 1).Value = 'Birmania'
$ws.Cells.Item(157, 3).Value = 4
$ws.Cells.Item(157, 4).Value = 0
$ws.Cells.Item(157, 5).Value = 14

$ws.Cells.Item(158, 1).Value = 'Bahamas'
$ws.Cells.Item(158, 2).Value = 14
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 1
$ws.Cells.Item(158, 5).Value = 13

$ws.Cells.Item(159, 1).Value = 'Dominica'
$ws.Cells.Item(159, 3).Value = 1

$ws.Cells.Item(170, 1).Value = 'Libia'

$ws.Cells.Item(172, 1).Value = 'Surinam'

$ws.Cells.Item(173, 1).Value = 'Mozambique'

$ws.Cells.Item(174, 1).Value = 'Laos'

$ws.Cells.Item(175, 1).Value = 'Seychelles'

$ws.Cells.Item(177, 1).Value = 'San Cristobal y Nieves'
$ws.Cells.Item(177, 3).Value = 5

$ws.Cells.Item(178, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(178, 3).Value = 0

$ws.Cells.Item(182, 1).Value = 'Santa Sede'

$ws.Cells.Item(183, 1).Value = 'Benin'

$ws.Cells.Item(184, 1).Value = 'San Martin (Parte Holandesa)'

$ws.Cells.Item(189, 1).Value = 'Republica del Chad'
$ws.Cells.Item(189, 3).Value = 2

$ws.Cells.Item(190, 1).Value = 'Fiyi'
$ws.Cells.Item(190, 3).Value = 0

$ws.Cells.Item(191, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(191, 3).Value = 1

$ws.Cells.Item(192, 1).Value = 'Montserrat'
$ws.Cells.Item(192, 3).Value = 0

$ws.Cells.Item(199, 1).Value = 'Republica de Africa Central'

$ws.Cells.Item(200, 1).Value = 'Somalia'

$ws.Cells.Item(201, 1).Value = 'Liberia'

$ws.Cells.Item(202, 1).Value = 'Anguila'

$ws.Cells.Item(203, 1).Value = 'Islas Virgenes Britanicas'
